# Update the "Keras Neural Network" best-score table cell on the
# "Algorithms - Results" slide from 54.60% to 56.08%.

$p = $ppt.ActivePresentation

$targetSlide = $null
$targetShape = $null
$targetRow = 0
$targetCol = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
                for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                    $cellText = $tbl.Cell($r, $c).Shape.TextFrame.TextRange.Text
                    if ($cellText -eq "54.60%") {
                        $targetSlide = $slide
                        $targetShape = $shape
                        $targetRow = $r
                        $targetCol = $c
                    }
                }
            }
        }
    }
}

if ($targetShape -ne $null) {
    $cell = $targetShape.Table.Cell($targetRow, $targetCol)
    $cell.Shape.TextFrame.TextRange.Text = "56.08%"
}
